$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of column letter -> new value (applies identically to rows 2 and 3)
$updates = @{
    "D"  = 0.163
    "G"  = -0.5034013605442178
    "H"  = -2.659863945578231
    "I"  = -3.034013605442177
    "J"  = -3.034013605442177
    "K"  = -6.37
    "L"  = -4.333333333333334
    "U"  = 0.26
    "V"  = 0.0143646408839779
    "W"  = 0.7406976744186047
    "X"  = 0.1421725538570773
    "Y"  = 0.5985251205615274
    "Z"  = 0.2094912355707567
    "AA" = -0.6355992589425679
    "AB" = 0.0838509377116391
    "AC" = -0.719450196654207
    "AD" = 18.5
    "AE" = 0
    "AF" = 18.5
    "AG" = 18.24
    "AH" = 0.5054644808743169
    "AI" = 1.161330822347772
    "AJ" = 0.5019262520638414
    "AK" = 1.164007657945118
    "AL" = 1.65
    "AM" = 1.65
    "AN" = -4.07488986784141
    "AO" = -2.703030303030303
    "AP" = -4.017621145374449
    "AQ" = -2.703030303030303
}

foreach ($col in $updates.Keys) {
    foreach ($row in 2..3) {
        $ws.Range("$col$row").Value = $updates[$col]
    }
}
